$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 191 (shifts old rows 191-219 down to 192-220)
$ws.Rows.Item(191).Insert()

# Fill in the new weekly price observation row
$ws.Range("A191").Value = 7
$ws.Range("B191").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C191").Value = "Ñuble"
$ws.Range("D191").Value = 44505
$ws.Range("E191").Value = 16
$ws.Range("F191").Value = 100114001
$ws.Range("G191").Value = "Papa"
$ws.Range("H191").Value = "Patagonia"
$ws.Range("I191").Value = "1a (guarda)"
$ws.Range("J191").Value = 260
$ws.Range("K191").Value = 7000
$ws.Range("L191").Value = 8000
$ws.Range("M191").Value = 7500
$ws.Range("N191").Value = "$/saco 25 kilos"
$ws.Range("O191").Value = "Provincia de Diguillín"
$ws.Range("P191").Value = 300
$ws.Range("Q191").Value = 25
$ws.Range("R191").Value = "Hortaliza"
